# Applies the data-refresh edit described in the commit:
# "idw pt.2 2018 - melhor com SUBestimacao a mesma - adicionar vies???? experimentar outros k e p"
# The whole data table (rows 2-47) is replaced with a new set of measurements
# dated 2018-06-20 (serial 43271) instead of 2017-10-18 (serial 43026), and two
# additional rows (48-49) are appended, extending the table to A1:C49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents for A2:C49 (Data, Estacao, pm10)
$rows = @(
    @(43271, 1054, 13.68181818181818),
    @(43271, 3093, 36.49166666666667),
    @(43271, 2017, 27.125),
    @(43271, 3075, 48.41666666666666),
    @(43271, 1053, 27.2),
    @(43271, 1052, 24.45454545454545),
    @(43271, 3104, 41.12380952380953),
    @(43271, 5012, 14.45833333333333),
    @(43271, 3096, 29.04583333333333),
    @(43271, 2006, 37.58333333333334),
    @(43271, 1046, 18.83333333333333),
    @(43271, 5011, 38.29166666666666),
    @(43271, 1048, 10.75),
    @(43271, 3072, 40.19583333333333),
    @(43271, 1023, 27.20833333333333),
    @(43271, 2019, 31.45833333333333),
    @(43271, 3095, 35.2875),
    @(43271, 2004, 17.89473684210526),
    @(43271, 3099, 38.3375),
    @(43271, 2021, 14.47368421052632),
    @(43271, 1028, 17.22727272727273),
    @(43271, 1042, 10.08333333333333),
    @(43271, 2020, 18.41666666666667),
    @(43271, 2018, 32.83333333333334),
    @(43271, 2016, 28.45833333333333),
    @(43271, 1030, 27.28571428571428),
    @(43271, 5007, 30.44166666666667),
    @(43271, 3083, 38.62916666666667),
    @(43271, 3055, 39.97916666666666),
    @(43271, 3085, 25.6375),
    @(43271, 3102, 25.84583333333333),
    @(43271, 5008, 43.41666666666666),
    @(43271, 1025, 20.25),
    @(43271, 3089, 27.72083333333333),
    @(43271, 1051, 16.58333333333333),
    @(43271, 2022, 28.625),
    @(43271, 3097, 43.25416666666666),
    @(43271, 3071, 36.34166666666667),
    @(43271, 1044, 27.79166666666667),
    @(43271, 3063, 38.65833333333333),
    @(43271, 1043, 7),
    @(43271, 3094, 39.85833333333333),
    @(43271, 3091, 28.62916666666667),
    @(43271, 3084, 31.62083333333333),
    @(43271, 1055, 25),
    @(43271, 4003, 19.29166666666667),
    @(43271, 4006, 25.33333333333333),
    @(43271, 1031, 25.70588235294118)
)

$firstDataRow = 2
$lastOriginalRow = 47

$r = $firstDataRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Rows beyond the original table (48, 49) need the date number format
    # explicitly applied to column A, matching the style used by the rest
    # of the column (existing rows already carry this style).
    if ($r -gt $lastOriginalRow) {
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }

    $r++
}
